$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 content updates (DungeonStory sheet, table row for id 47000021 "远古之路")
$ws.Range("I8").Value = "bossmanwang=bosstalic,potteryroom=suntemple,trapspear=trapspearwall,stonedoor2=snowhill"
$ws.Range("M8").Value = "fightfail"
$ws.Range("M8").Font.Color = 255
$ws.Range("C8").Value = "从冰冻苔原往前走,在亚瑞特山脚下，就是传说中的远古之路。再往前就是亚瑞特山脉颠峰，在那里有3个古代勇士守护者通往远古遗迹的道路。|n#ff6666|你需要激活机关，并打败3个野蛮人的灵魂|n|n#cccc66|★关卡中事件难度更大|n#66cccc|★开始附带诅咒【胆小鬼】"

# Move the active selection to C8, matching the saved view state
$ws.Range("C8").Select()
